# Global controls partially implemented.
# - "current punchlist" sheet: row 15 gains a DONE / date / version entry for
#   the "Move to github." item, and a new row 18 documents a follow-up
#   priority-9 item: "Deploy from github directly."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 already has A15 (priority 4) and E15 ("Move to github."). Fill in
# the status/date/version columns for that item.
$ws.Range("B15").Value = "DONE"

# Give C15 the same date style already used by the other date cells in the
# column (xfId 3 / numFmtId 14) before writing the serial value, so a brand
# new cell doesn't pick up an auto-detected number format.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value = 44671

$ws.Range("D15").Value = "0.9.1"

# New row 18: a follow-up task.
$ws.Range("A18").Value = 9
$ws.Range("E18").Value = "Deploy from github directly."

# Match the saved selection/cursor position.
$ws.Range("E18").Select() | Out-Null
